# Add two new columns "I0" (column I) and "IF" (column J) to the sheet,
# mirroring the header style used by the existing "IP" header (column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, matching the style of the existing header cell H1
# (bold font, thin border, centered alignment). Copy/PasteSpecial reuses
# the existing style record instead of minting a new (slightly different)
# one the way direct Font/Border property assignment would.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-26 (column I = I0, column J = IF).
$data = @(
    @(1, 2),
    @(8, 8),
    @(5, 6),
    @(4, 5),
    @(2, 3),
    @(5, 5),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(8, 9),
    @(4, 4),
    @(7, 8),
    @(6, 8),
    @(4, 5),
    @(5, 5),
    @(5, 5),
    @(5, 6),
    @(4, 5),
    @(8, 9),
    @(2, 2),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
